{"js": "// Disposals approval template: the \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.:\" (protocol number) line\n// referenced the merge field ${regionaldirect_protocol} directly after the\n// \".: \" label. The edit adds a literal protocol-type prefix \"\u03a6.11.2/\"\n// immediately before that merge field, so the line now reads\n// \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.11.2/${regionaldirect_protocol}\" instead of\n// \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ${regionaldirect_protocol}\".\n\nconst body = context.document.body;\n\n// \"${regionaldirect_protocol}\" is unique in the document (the only other\n// similar placeholder is \"${regionaldirect_protocoldate}\" a few lines above,\n// which is untouched by this change because it ends in \"date}\" not \"}\").\nconst results = body.search(\"${regionaldirect_protocol}\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the ${regionaldirect_protocol} placeholder\");\n}\n\n// Insert the new literal text right before the placeholder, preserving the\n// placeholder range itself (and everything else in the paragraph).\nconst target = results.items[0];\ntarget.insertText(\"\u03a6.11.2/\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Disposals approval template: the \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.:\" (protocol number) line\n# referenced the merge field ${regionaldirect_protocol} directly after the\n# \".: \" label. This edit adds a literal protocol-type prefix \"\u03a6.11.2/\"\n# immediately before that merge field, so the line now reads\n# \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.11.2/${regionaldirect_protocol}\" instead of\n# \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ${regionaldirect_protocol}\".\n\n$d = $word.ActiveDocument\n\n# \"${regionaldirect_protocol}\" is unique in the document (the only other\n# similar placeholder is \"${regionaldirect_protocoldate}\" a few lines above,\n# which is untouched by this change because it ends in \"date}\", not \"}\").\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$searchText  = \"`${regionaldirect_protocol}\"\n$replaceText = \"\u03a6.11.2/`${regionaldirect_protocol}\"\n\n$wdFindContinue = 1\n$wdReplaceOne   = 1\n\n$result = $find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne)\n\nif (-not $result) {\n    throw \"Could not find the `${regionaldirect_protocol} placeholder\"\n}\n"}
